# Databases.xlsx -- add the "Card Drops" sheet (card-generation weights table)
# after the existing "Enemies" sheet, and make it the active tab.

$wb = $excel.ActiveWorkbook
$enemies = $wb.Worksheets.Item("Enemies")

# New sheet goes right after "Enemies" -> sheetId 2, second tab.
$cards = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $enemies)
$cards.Name = "Card Drops"

# --- Labels (column A) ---------------------------------------------------
$cards.Range("A1").Value = "Main Spec"
$cards.Range("A2").Value = "Off Spec"
$cards.Range("A3").Value = "Potions"
$cards.Range("A4").Value = "Other Specs"
$cards.Range("A5").Value = "Relics"
$cards.Range("A7").Value = "Total"

# --- Weights (column B), formatted as a percentage ------------------------
$cards.Range("B1").Value = 0.4
$cards.Range("B2").Value = 0.2
$cards.Range("B3").Value = 0.2
$cards.Range("B4").Value = 0.1
$cards.Range("B5").Value = 0.1
$cards.Range("B7").Formula = "=SUM(B1:B5)"
$cards.Range("B1:B7").NumberFormat = "0%"

# --- Note next to the Relics row ------------------------------------------
$cards.Range("D5").Value = "60/30/10% Main/Off/Other Stat Relic"

# --- Decorative legend / swatch block (I1:L8) ------------------------------
# Cells that carry a white font (used over a colored swatch fill).
$cards.Range("I1:K2").Font.ThemeColor = 2
$cards.Range("I4:K4").Font.ThemeColor = 2
$cards.Range("J3").Font.ThemeColor = 2

# Remaining cells in the block are present but carry no special formatting;
# touch a no-op border property so the (blank) cell still gets emitted.
# (done one column at a time -- applying it to a multi-column range only
# forces serialization of the first column)
$cards.Range("L1:L8").Borders.Item(7).LineStyle = -4142
$cards.Range("I3").Borders.Item(7).LineStyle = -4142
$cards.Range("K3").Borders.Item(7).LineStyle = -4142
$cards.Range("I5:I8").Borders.Item(7).LineStyle = -4142
$cards.Range("J5:J8").Borders.Item(7).LineStyle = -4142
$cards.Range("K5:K8").Borders.Item(7).LineStyle = -4142

# --- Sheet/view cosmetics ---------------------------------------------------
# (target stored width is 13.140625 char-units; the engine quantizes
# ColumnWidth input to steps of 1/6, so 12.3 is the closest reachable input)
$cards.Columns.Item(1).ColumnWidth = 12.3
$cards.PageSetup.Orientation = 1
[void]$cards.Range("H6").Select()

# Card Drops becomes the active tab; Enemies loses its "tabSelected" flag.
$cards.Activate()

Write-Host "Card Drops sheet added"
